$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.171.53"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "2.790.64"
$ws.Range("E3").Value = "  +5.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'116.94"
$ws.Range("E5").Value = "  +4.49%  "
$ws.Range("D6").Value = "'341.36"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("D7").Value = "'0.555"
$ws.Range("E7").Value = "  +5.82%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  +6.03%  "
$ws.Range("D10").Value = "'42.10"
$ws.Range("E10").Value = "  +6.37%  "
$ws.Range("E11").Value = "  +7.33%  "
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").Value = "'7.62"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "3.229.79"
$ws.Range("E15").Value = "  +5.79%  "
$ws.Range("D16").Value = "2.804.58"
$ws.Range("E16").Value = "  +6.38%  "
$ws.Range("D17").Value = "'0.886"
$ws.Range("E17").Value = "  +3.82%  "
$ws.Range("D18").Value = "52.032.73"
$ws.Range("E18").Value = "  +5.04%  "
$ws.Range("D19").Value = "'3.22"
$ws.Range("E19").Value = "  +10.54%  "
$ws.Range("D20").Value = "'13.35"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("D22").Value = "0.0₃0987"
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").Value = "'278.18"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").Value = "'70.34"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "'2.81"
$ws.Range("E25").Value = "  +10.03%  "
$ws.Range("D26").Value = "'26.86"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'10.18"
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D31").Value = "'34.88"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").Value = "'50.53"
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("D33").Value = "'5.71"
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("D34").Value = "'0.0826"
$ws.Range("E34").Value = "  +2.29%  "
$ws.Range("E35").Value = "  +4.70%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'18.99"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'4.97"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'3.27"
$ws.Range("E39").Value = "  +5.76%  "
$ws.Range("E40").Value = "  +28.76%  "
$ws.Range("D41").Value = "'0.0374"
$ws.Range("E41").Value = "  +13.77%  "
$ws.Range("E42").Value = "  +5.13%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.116"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'126.88"
$ws.Range("E44").Value = "  -1.60%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'23.11"
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("D46").Value = "2.105.22"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D49").Value = "'5.55"
$ws.Range("E49").Value = "  +6.36%  "
$ws.Range("D50").Value = "'0.911"
$ws.Range("E50").Value = "  +21.01%  "
$ws.Range("D51").Value = "'8.95"
$ws.Range("E51").Value = "  +0.95%  "
